$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.364.49'
$ws.Range("E2").Value = '  -3.82%  '
$ws.Range("D3").Value = '3.304.36'
$ws.Range("E3").Value = '  -4.13%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.01'
$ws.Range("E5").Value = '  -3.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.75'
$ws.Range("E6").Value = '  -4.02%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.306.24'
$ws.Range("E8").Value = '  -4.11%  '
$ws.Range("E9").Value = '  -2.20%  '
$ws.Range("E10").Value = '  -2.63%  '
$ws.Range("E11").Value = '  -3.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.408'
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("D13").Value = '3.855.25'
$ws.Range("E13").Value = '  -4.61%  '
$ws.Range("E14").Value = '  +0.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.19'
$ws.Range("E15").Value = '  -4.04%  '
$ws.Range("D16").Value = '3.303.25'
$ws.Range("E16").Value = '  -5.46%  '
$ws.Range("E17").Value = '  -3.78%  '
$ws.Range("D18").Value = '60.221.57'
$ws.Range("E18").Value = '  -4.19%  '
$ws.Range("E19").Value = '  -5.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.31'
$ws.Range("E20").Value = '  -1.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.57'
$ws.Range("E21").Value = '  -4.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '375.83'
$ws.Range("E22").Value = '  -3.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.82'
$ws.Range("E23").Value = '  -1.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.545'
$ws.Range("E24").Value = '  -4.22%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  -4.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000104'
$ws.Range("E27").Value = '  -9.17%  '
$ws.Range("E28").Value = '  -6.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.24'
$ws.Range("E30").Value = '  -5.96%  '
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.66'
$ws.Range("E32").Value = '  -4.40%  '
$ws.Range("E33").Value = '  -4.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '22.55'
$ws.Range("E34").Value = '  -3.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.26'
$ws.Range("E35").Value = '  -6.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.23'
$ws.Range("E36").Value = '  -3.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '166.56'
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.74'
$ws.Range("E38").Value = '  -2.92%  '
$ws.Range("E39").Value = '  -7.72%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.92'
$ws.Range("E40").Value = '  -15.74%  '
$ws.Range("E41").Value = '  -4.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0742'
$ws.Range("E42").Value = '  -5.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.83'
$ws.Range("E43").Value = '  -2.14%  '
$ws.Range("E44").Value = '  -4.41%  '
$ws.Range("E45").Value = '  -4.40%  '
$ws.Range("E46").Value = '  -5.89%  '
$ws.Range("E47").Value = '  -4.56%  '
$ws.Range("D48").Value = '2.353.39'
$ws.Range("E48").Value = '  -7.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.998'
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("E50").Value = '  -5.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0256'
$ws.Range("E51").Value = '  -4.00%  '
